$wb = $excel.ActiveWorkbook

# --- Layer0 sheet: grid B2:Q9 (Bias/Weight0-7 rows x Neuron0-15 columns) ---
$ws1 = $wb.Worksheets.Item("Layer0")

$grid1 = New-Object 'object[,]' 8,16
$grid1[0,0] = -80.43898936921897
$grid1[0,1] = -84.59693280944545
$grid1[0,2] = 21.34236523684379
$grid1[0,3] = 10.38196275134003
$grid1[0,4] = -95.12887749917114
$grid1[0,5] = -92.19922354081878
$grid1[0,6] = 26.32510524637755
$grid1[0,7] = -87.49995246128725
$grid1[0,8] = -36.51126466649291
$grid1[0,9] = 41.22878403531359
$grid1[0,10] = -66.39558306760348
$grid1[0,11] = -35.82270513158981
$grid1[0,12] = 6.547334700058046
$grid1[0,13] = -57.86166776779288
$grid1[0,14] = -5.303354547393592
$grid1[0,15] = 25.49832450630915
$grid1[1,0] = 36.12173736708888
$grid1[1,1] = 36.6965667762617
$grid1[1,2] = -23.44653238773396
$grid1[1,3] = 15.21481609898233
$grid1[1,4] = 22.64823500236385
$grid1[1,5] = 31.17960714002092
$grid1[1,6] = 37.55553220181722
$grid1[1,7] = 24.69445985098462
$grid1[1,8] = 37.65414465796103
$grid1[1,9] = 34.72733634546212
$grid1[1,10] = 40.08828477343197
$grid1[1,11] = 50.9459639680981
$grid1[1,12] = 52.17678629191545
$grid1[1,13] = 48.27450336493471
$grid1[1,14] = 48.04927139684878
$grid1[1,15] = 24.83091110107623
$grid1[2,0] = 31.93669440890826
$grid1[2,1] = 27.42958055239869
$grid1[2,2] = 22.98535237417797
$grid1[2,3] = 5.474406718883306
$grid1[2,4] = 23.76700586937886
$grid1[2,5] = 27.11280669802687
$grid1[2,6] = -0.5426471678387098
$grid1[2,7] = 23.87075216396777
$grid1[2,8] = 12.55619769111444
$grid1[2,9] = -3.92277470880153
$grid1[2,10] = 30.63898971081599
$grid1[2,11] = 43.52161071855134
$grid1[2,12] = 1.316285222177179
$grid1[2,13] = 40.02004666684572
$grid1[2,14] = 40.75862693202311
$grid1[2,15] = 38.31820406273182
$grid1[3,0] = 10.84838020883303
$grid1[3,1] = 18.39450958226062
$grid1[3,2] = 5.964016709850147
$grid1[3,3] = 23.9704887004397
$grid1[3,4] = 16.49886809624139
$grid1[3,5] = 16.20013579152447
$grid1[3,6] = 2.767597077344891
$grid1[3,7] = 12.35497611026126
$grid1[3,8] = 26.56690539195412
$grid1[3,9] = 1.480528378347572
$grid1[3,10] = 14.47997330266823
$grid1[3,11] = 14.193358526084
$grid1[3,12] = 9.232601559639386
$grid1[3,13] = 13.45533195688077
$grid1[3,14] = 4.506013323324806
$grid1[3,15] = 0.9495051997838768
$grid1[4,0] = 2.40863374808568
$grid1[4,1] = 16.53777520442131
$grid1[4,2] = 0.7189910774185907
$grid1[4,3] = 14.58348644402134
$grid1[4,4] = 8.721856521101261
$grid1[4,5] = 16.42272544058079
$grid1[4,6] = 2.485669232755435
$grid1[4,7] = 12.21574172971754
$grid1[4,8] = 19.31443802007323
$grid1[4,9] = 1.495710335843043
$grid1[4,10] = 10.79192134704409
$grid1[4,11] = 8.638529127920561
$grid1[4,12] = 6.238474256616009
$grid1[4,13] = 7.889235635674853
$grid1[4,14] = 6.787766561437699
$grid1[4,15] = -4.921624540120979
$grid1[5,0] = 1.386852788308471
$grid1[5,1] = 1.674038817980375
$grid1[5,2] = 3.310216006770275
$grid1[5,3] = 4.819300480219698
$grid1[5,4] = 0.9248943133802958
$grid1[5,5] = 2.024029031774504
$grid1[5,6] = -0.1502223852548318
$grid1[5,7] = 2.83263571191881
$grid1[5,8] = 3.296912335015902
$grid1[5,9] = 0.03610241284409525
$grid1[5,10] = 2.249625007395552
$grid1[5,11] = 2.414590305445437
$grid1[5,12] = -1.382018293644014
$grid1[5,13] = 1.467238525330903
$grid1[5,14] = 0.5226543081061649
$grid1[5,15] = 0.3310723279829608
$grid1[6,0] = 0.234886022788448
$grid1[6,1] = 2.307973429121319
$grid1[6,2] = -1.30811233217117
$grid1[6,3] = 0.5616666195211862
$grid1[6,4] = 1.965838236257432
$grid1[6,5] = 1.985181297714537
$grid1[6,6] = 0.69750481803818
$grid1[6,7] = -0.487758710887738
$grid1[6,8] = 4.688530991267798
$grid1[6,9] = -0.123269317262097
$grid1[6,10] = 2.13186204198269
$grid1[6,11] = -1.003966738577074
$grid1[6,12] = 2.106026688975086
$grid1[6,13] = 0.8787126322624352
$grid1[6,14] = -0.6086727745488361
$grid1[6,15] = 1.844944573837262
$grid1[7,0] = 0.5875991551779893
$grid1[7,1] = 1.229408458289134
$grid1[7,2] = -0.6809125198675811
$grid1[7,3] = -1.046585594789087
$grid1[7,4] = 1.297238711683581
$grid1[7,5] = 1.929372806937975
$grid1[7,6] = -0.1446273075641309
$grid1[7,7] = 2.026026573466379
$grid1[7,8] = -3.200875860748735
$grid1[7,9] = 0.5468405783365957
$grid1[7,10] = 1.338712719398873
$grid1[7,11] = 1.96868303023525
$grid1[7,12] = 1.292456659416576
$grid1[7,13] = 4.18929850050441
$grid1[7,14] = 3.080578912428321
$grid1[7,15] = 1.66130742878565

$ws1.Range("B2:Q9").Value = $grid1

# --- Layer1 sheet: single column B2:B18 (Bias/Weight0-15) ---
$ws2 = $wb.Worksheets.Item("Layer1")

$grid2 = New-Object 'object[,]' 17,1
$grid2[0,0] = -130.5267979698347
$grid2[1,0] = 237.3162127402424
$grid2[2,0] = 241.0222913309289
$grid2[3,0] = 198.2821385838752
$grid2[4,0] = 148.1299780180685
$grid2[5,0] = 255.1211460002976
$grid2[6,0] = 253.0847800494708
$grid2[7,0] = 206.8081961141725
$grid2[8,0] = 239.7413955284346
$grid2[9,0] = 189.5307542810294
$grid2[10,0] = 217.9731197395053
$grid2[11,0] = 237.0091498950191
$grid2[12,0] = 219.0983079764121
$grid2[13,0] = 187.4066649630552
$grid2[14,0] = 236.6716639065425
$grid2[15,0] = 204.0489046305277
$grid2[16,0] = 183.5250503332178

$ws2.Range("B2:B18").Value = $grid2

